# Add the missing GC field definitions to "gc_fields_display" and
# re-point the active sheet / selections as left by the author after
# reviewing the new rows.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("gc_fields_display")
$ws2 = $wb.Worksheets.Item("gc_fields_uom")

# field (column A) / display name (column B) pairs, appended after the
# existing data (which currently ends at row 52).
$fieldNames = @(
    "GainCorrectedElevation",
    "GainUncorrectedElevation",
    "LossCorrectedElevation",
    "LossUncorrectedElevation",
    "MaxCorrectedElevation",
    "MaxGroundContactTime",
    "MaxMomentaryEnergyExpenditure",
    "MaxRelativeRunningEconomy",
    "MaxUncorrectedElevation",
    "MaxVerticalOscillation",
    "MaxVerticalSpeed",
    "MinBikeCadence",
    "MinCorrectedElevation",
    "MinGroundContactTime",
    "MinMomentaryEnergyExpenditure",
    "MinPace",
    "MinRelativeRunningEconomy",
    "MinRunCadence",
    "MinUncorrectedElevation",
    "MinVerticalOscillation",
    "SumAnaerobicTrainingEffect",
    "SumFloorsClimbed",
    "SumFloorsDescended",
    "WeightedMeanLeftBalance",
    "WeightedMeanMomentaryEnergyExpenditure",
    "WeightedMeanRelativeRunningEconomy",
    "WeightedMeanStanceTime",
    "WeightedMeanStanceTimePercent",
    "WeightedMeanVerticalSpeed"
)

$displayNames = @(
    "Elevation Gain",
    "Elevation Gain",
    "Elevation Loss",
    "Elevation Loss",
    "Max Elevation",
    "Max Ground Contact Time",
    "Max Energy Expenditure",
    "Max Running Economy",
    "Max Elevation",
    "Max Vertical Oscillation",
    "Max Vertical Speed",
    "Min Cadence",
    "Min Elevation",
    "Min Ground Contact Time",
    "Min Energy Expenditure",
    "Min Pace",
    "Min Running Economy",
    "Min Cadence",
    "Min Elevation",
    "Min Vertical Oscillation",
    "Anaerobic Training Effect",
    "Floors Climbed",
    "Floors Descended",
    "Avg Left Balance",
    "Avg Energy Expenditure",
    "Avg Running Economy",
    "Avg Stance Time",
    "Avg Stance Percent",
    "Avg Vertical Speed"
)

$startRow = 53

# Write column A (field) fully first ...
$r = $startRow
foreach ($name in $fieldNames) {
    $ws1.Cells.Item($r, 1).Value = $name
    $r++
}

# ... then column B (display name), matching the order the rows were
# originally authored in.
$r = $startRow
foreach ($name in $displayNames) {
    $ws1.Cells.Item($r, 2).Value = $name
    $r++
}

# Leave the workbook as the author did: "gc_fields_uom" selection moved
# to C20 (no longer the active tab) and "gc_fields_display" becomes the
# active tab with the new last row selected.
$ws2.Activate()
$ws2.Range("C20").Select()

$ws1.Activate()
$ws1.Range("B82").Select()
